$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = 45371
$ws.Range("F3").Value = 45370
$ws.Range("F4").Value = 45369
$ws.Range("F5").Value = 45368
$ws.Range("F6").Value = 45367
$ws.Range("F7").Value = 45366
